$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for B:G across rows 2-6 (regenerated s_vals after filtering save games)
$data = @{
    2 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 0, 8.418600821238126)
    3 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 6.048734245549538)
    4 = @(0.1554434735375247, 1.65323645889881, 0.7127328510149897, 6.48142807727062, 1, 9.002840860721944)
    5 = @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 4.371470058157054)
    6 = @(0.0003714022599530242, 0.05231270169004087, 16.98373111632243, 0.4998867070740569, 1, 17.53630192734648)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = $i + 2  # B = 2
        $ws.Cells.Item($row, $col).Value = $vals[$i]
    }
}
